# Regen sval data to filter save games
# Update B2:E5 (and derived G2:G5 "sum" column) with new values per the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.230985683306322;    C = 1.667794583268128;   D = 0.8054896365839992;   E = 0.496779210170732 }
    3 = @{ B = 3.230985683306322;    C = 1.667794583268128;   D = 0.8054896365839992;   E = 0.496779210170732 }
    4 = @{ B = 0.0008583669626518464; C = 0.3127903958511391; D = 3.900430680208489;    E = 8.660232485948974 }
    5 = @{ B = 0.6753301551942219;   C = 1.667794583268128;   D = 3.900430680208489;    E = 0.496779210170732 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.B + $vals.C + $vals.D + $vals.E
}
